$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The wording of the fourth test case (row 5) was tweaked: remove the
# "on the Options part" phrase so it reads "... for issues from the
# Issues tab and applying the changes".
$ws.Range("B5").Value = 'Adding an "Author" column for issues from the "Issues" tab and applying the changes'

# The view was scrolled down a bit and a different cell became selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("C9").Select()
